# Append 12 new Google-Forms survey responses (rows 233-244) to the single
# worksheet "Form Responses 1", matching new submissions received after the
# workbook was last exported. Commit message: "hacer descargables las bases".
#
# Each new row reuses the formatting already present on an existing, fully
# populated data row (row 4): column A keeps the date/time number format,
# every other column keeps the plain text/number format. Per-cell Copy +
# PasteSpecial(xlPasteFormats) is used (rather than pasting an entire A:AC
# block at once) so that columns left blank in a given response do not end
# up as empty-but-styled cells - they simply stay absent, exactly as in the
# other rows of this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormatSource = $ws.Range("A4")   # column A: timestamp number format
$textFormatSource  = $ws.Range("B4")   # every other column: plain format
$xlPasteFormats = -4122

# ----- Row 233 -----
$dateFormatSource.Copy()
$ws.Range("A233").PasteSpecial($xlPasteFormats)
$ws.Range("A233").Value = 44697.486502569445
$textFormatSource.Copy()
$ws.Range("B233").PasteSpecial($xlPasteFormats)
$ws.Range("B233").Value = "carolina.timana.cfk2022@gmail.com"
$textFormatSource.Copy()
$ws.Range("C233").PasteSpecial($xlPasteFormats)
$ws.Range("C233").Value = 140.0
$textFormatSource.Copy()
$ws.Range("D233").PasteSpecial($xlPasteFormats)
$ws.Range("D233").Value = "Una"
$textFormatSource.Copy()
$ws.Range("E233").PasteSpecial($xlPasteFormats)
$ws.Range("E233").Value = "No"
$textFormatSource.Copy()
$ws.Range("F233").PasteSpecial($xlPasteFormats)
$ws.Range("F233").Value = 0.0
$textFormatSource.Copy()
$ws.Range("G233").PasteSpecial($xlPasteFormats)
$ws.Range("G233").Value = "Televisor, Video proyector"
$textFormatSource.Copy()
$ws.Range("H233").PasteSpecial($xlPasteFormats)
$ws.Range("H233").Value = "Redes eléctricas"
$textFormatSource.Copy()
$ws.Range("I233").PasteSpecial($xlPasteFormats)
$ws.Range("I233").Value = "Existe una red de internet de claro en  el programa centros digitales pero no funciona se quemo el receptor y la estructura esta dañada"
$textFormatSource.Copy()
$ws.Range("J233").PasteSpecial($xlPasteFormats)
$ws.Range("J233").Value = 122.0
$textFormatSource.Copy()
$ws.Range("K233").PasteSpecial($xlPasteFormats)
$ws.Range("K233").Value = 122.0
$textFormatSource.Copy()
$ws.Range("L233").PasteSpecial($xlPasteFormats)
$ws.Range("L233").Value = 1.0
$textFormatSource.Copy()
$ws.Range("M233").PasteSpecial($xlPasteFormats)
$ws.Range("M233").Value = "No"
$textFormatSource.Copy()
$ws.Range("N233").PasteSpecial($xlPasteFormats)
$ws.Range("N233").Value = "Windows 7 o superior (64 bits)"
$textFormatSource.Copy()
$ws.Range("O233").PasteSpecial($xlPasteFormats)
$ws.Range("O233").Value = 0.0
$textFormatSource.Copy()
$ws.Range("P233").PasteSpecial($xlPasteFormats)
$ws.Range("P233").Value = "Paquete de oficina (Ej. Microsoft Office), Lenguaje de programación (Ej. MakeCode), Paquete de diseño gráfico (Ej. Programas de edición de video, edición de mapas de bits, etc), Paquete de software educativo (Ej. cursos de inglés)"
$textFormatSource.Copy()
$ws.Range("Q233").PasteSpecial($xlPasteFormats)
$ws.Range("Q233").Value = "El internet es necesario para esta institución y el no tenerlo les limitan muchas capacitaciones a nivel de formación docente cómo para la enseñanza a sus estudiantes"
$textFormatSource.Copy()
$ws.Range("R233").PasteSpecial($xlPasteFormats)
$ws.Range("R233").Value = "Sí"
$textFormatSource.Copy()
$ws.Range("S233").PasteSpecial($xlPasteFormats)
$ws.Range("S233").Value = 24.0
$textFormatSource.Copy()
$ws.Range("T233").PasteSpecial($xlPasteFormats)
$ws.Range("T233").Value = "No"
$textFormatSource.Copy()
$ws.Range("U233").PasteSpecial($xlPasteFormats)
$ws.Range("U233").Value = 0.0
$textFormatSource.Copy()
$ws.Range("V233").PasteSpecial($xlPasteFormats)
$ws.Range("V233").Value = "Tarjetas micro:bit."
$textFormatSource.Copy()
$ws.Range("W233").PasteSpecial($xlPasteFormats)
$ws.Range("W233").Value = 4.0
$textFormatSource.Copy()
$ws.Range("X233").PasteSpecial($xlPasteFormats)
$ws.Range("X233").Value = "No hay ningún elemento adicional a las micro:bits"
$textFormatSource.Copy()
$ws.Range("Y233").PasteSpecial($xlPasteFormats)
$ws.Range("Y233").Value = "Existen muchos equipos portátiles pero no existen instalaciones adecuadas para utilizarlos partiendo que la energía con la que se alimenta la institución es solar faltarían más paneles para proveer una mejor fuente de energía "
$textFormatSource.Copy()
$ws.Range("Z233").PasteSpecial($xlPasteFormats)
$ws.Range("Z233").Value = "No"
$textFormatSource.Copy()
$ws.Range("AA233").PasteSpecial($xlPasteFormats)
$ws.Range("AA233").Value = "No"
$textFormatSource.Copy()
$ws.Range("AB233").PasteSpecial($xlPasteFormats)
$ws.Range("AB233").Value = "No"

# ----- Row 234 -----
$dateFormatSource.Copy()
$ws.Range("A234").PasteSpecial($xlPasteFormats)
$ws.Range("A234").Value = 44697.49317667824
$textFormatSource.Copy()
$ws.Range("B234").PasteSpecial($xlPasteFormats)
$ws.Range("B234").Value = "anabell.zuniga.cfk2022@gmail.com"
$textFormatSource.Copy()
$ws.Range("C234").PasteSpecial($xlPasteFormats)
$ws.Range("C234").Value = 24.0
$textFormatSource.Copy()
$ws.Range("D234").PasteSpecial($xlPasteFormats)
$ws.Range("D234").Value = "Una"
$textFormatSource.Copy()
$ws.Range("E234").PasteSpecial($xlPasteFormats)
$ws.Range("E234").Value = "No"
$textFormatSource.Copy()
$ws.Range("F234").PasteSpecial($xlPasteFormats)
$ws.Range("F234").Value = 0.0
$textFormatSource.Copy()
$ws.Range("G234").PasteSpecial($xlPasteFormats)
$ws.Range("G234").Value = "No hay más dispositivos en las aulas"
$textFormatSource.Copy()
$ws.Range("H234").PasteSpecial($xlPasteFormats)
$ws.Range("H234").Value = "Redes eléctricas"
$textFormatSource.Copy()
$ws.Range("I234").PasteSpecial($xlPasteFormats)
$ws.Range("I234").Value = "No cuentan con infraestructura de Telecomunicaciones"
$textFormatSource.Copy()
$ws.Range("J234").PasteSpecial($xlPasteFormats)
$ws.Range("J234").Value = 55.0
$textFormatSource.Copy()
$ws.Range("K234").PasteSpecial($xlPasteFormats)
$ws.Range("K234").Value = 55.0
$textFormatSource.Copy()
$ws.Range("L234").PasteSpecial($xlPasteFormats)
$ws.Range("L234").Value = 0.0
$textFormatSource.Copy()
$ws.Range("M234").PasteSpecial($xlPasteFormats)
$ws.Range("M234").Value = "Sí"
$textFormatSource.Copy()
$ws.Range("N234").PasteSpecial($xlPasteFormats)
$ws.Range("N234").Value = "Windows 7 o superior (32 bits), Windows 7 o superior (64 bits)"
$textFormatSource.Copy()
$ws.Range("O234").PasteSpecial($xlPasteFormats)
$ws.Range("O234").Value = 0.0
$textFormatSource.Copy()
$ws.Range("P234").PasteSpecial($xlPasteFormats)
$ws.Range("P234").Value = "Paquete de oficina (Ej. Microsoft Office), Lenguaje de programación (Ej. MakeCode), Paquete de software educativo (Ej. cursos de inglés)"
$textFormatSource.Copy()
$ws.Range("R234").PasteSpecial($xlPasteFormats)
$ws.Range("R234").Value = "Sí"
$textFormatSource.Copy()
$ws.Range("S234").PasteSpecial($xlPasteFormats)
$ws.Range("S234").Value = 22.0
$textFormatSource.Copy()
$ws.Range("T234").PasteSpecial($xlPasteFormats)
$ws.Range("T234").Value = "Sí"
$textFormatSource.Copy()
$ws.Range("U234").PasteSpecial($xlPasteFormats)
$ws.Range("U234").Value = 0.0
$textFormatSource.Copy()
$ws.Range("V234").PasteSpecial($xlPasteFormats)
$ws.Range("V234").Value = "Tarjetas micro:bit."
$textFormatSource.Copy()
$ws.Range("W234").PasteSpecial($xlPasteFormats)
$ws.Range("W234").Value = 12.0
$textFormatSource.Copy()
$ws.Range("X234").PasteSpecial($xlPasteFormats)
$ws.Range("X234").Value = "Sensores, Soldadura electrónica, Cautines, Pelacables, Pinzas"
$textFormatSource.Copy()
$ws.Range("Y234").PasteSpecial($xlPasteFormats)
$ws.Range("Y234").Value = "Disponen de un kit STEM"
$textFormatSource.Copy()
$ws.Range("Z234").PasteSpecial($xlPasteFormats)
$ws.Range("Z234").Value = "No"
$textFormatSource.Copy()
$ws.Range("AA234").PasteSpecial($xlPasteFormats)
$ws.Range("AA234").Value = "Sí"
$textFormatSource.Copy()
$ws.Range("AB234").PasteSpecial($xlPasteFormats)
$ws.Range("AB234").Value = "No"
$textFormatSource.Copy()
$ws.Range("AC234").PasteSpecial($xlPasteFormats)
$ws.Range("AC234").Value = "No tienen kits desconectados. Cuentan con 1 Kit STEM y las microbits que obtuvieron los docentes, formados en vigencias anteriores."

# ----- Row 235 -----
$dateFormatSource.Copy()
$ws.Range("A235").PasteSpecial($xlPasteFormats)
$ws.Range("A235").Value = 44697.53576863426
$textFormatSource.Copy()
$ws.Range("B235").PasteSpecial($xlPasteFormats)
$ws.Range("B235").Value = "camilo.villota.cfk2022@gmail.com"
$textFormatSource.Copy()
$ws.Range("C235").PasteSpecial($xlPasteFormats)
$ws.Range("C235").Value = 40.0
$textFormatSource.Copy()
$ws.Range("D235").PasteSpecial($xlPasteFormats)
$ws.Range("D235").Value = "Tres"
$textFormatSource.Copy()
$ws.Range("E235").PasteSpecial($xlPasteFormats)
$ws.Range("E235").Value = "No"
$textFormatSource.Copy()
$ws.Range("F235").PasteSpecial($xlPasteFormats)
$ws.Range("F235").Value = 0.0
$textFormatSource.Copy()
$ws.Range("G235").PasteSpecial($xlPasteFormats)
$ws.Range("G235").Value = "Televisor"
$textFormatSource.Copy()
$ws.Range("H235").PasteSpecial($xlPasteFormats)
$ws.Range("H235").Value = "Redes eléctricas, Puertos de cableado estructurado para conexión a la red"
$textFormatSource.Copy()
$ws.Range("I235").PasteSpecial($xlPasteFormats)
$ws.Range("I235").Value = "No se cuenta con el servicio de Internet"
$textFormatSource.Copy()
$ws.Range("J235").PasteSpecial($xlPasteFormats)
$ws.Range("J235").Value = 60.0
$textFormatSource.Copy()
$ws.Range("K235").PasteSpecial($xlPasteFormats)
$ws.Range("K235").Value = 60.0
$textFormatSource.Copy()
$ws.Range("L235").PasteSpecial($xlPasteFormats)
$ws.Range("L235").Value = 0.0
$textFormatSource.Copy()
$ws.Range("M235").PasteSpecial($xlPasteFormats)
$ws.Range("M235").Value = "No"
$textFormatSource.Copy()
$ws.Range("N235").PasteSpecial($xlPasteFormats)
$ws.Range("N235").Value = "Windows 7 o superior (64 bits)"
$textFormatSource.Copy()
$ws.Range("O235").PasteSpecial($xlPasteFormats)
$ws.Range("O235").Value = 0.0
$textFormatSource.Copy()
$ws.Range("P235").PasteSpecial($xlPasteFormats)
$ws.Range("P235").Value = "Paquete de oficina (Ej. Microsoft Office), Lenguaje de programación (Ej. MakeCode), Paquete de diseño gráfico (Ej. Programas de edición de video, edición de mapas de bits, etc)"
$textFormatSource.Copy()
$ws.Range("R235").PasteSpecial($xlPasteFormats)
$ws.Range("R235").Value = "Sí"
$textFormatSource.Copy()
$ws.Range("S235").PasteSpecial($xlPasteFormats)
$ws.Range("S235").Value = 60.0
$textFormatSource.Copy()
$ws.Range("T235").PasteSpecial($xlPasteFormats)
$ws.Range("T235").Value = "Sí"
$textFormatSource.Copy()
$ws.Range("U235").PasteSpecial($xlPasteFormats)
$ws.Range("U235").Value = 0.0
$textFormatSource.Copy()
$ws.Range("V235").PasteSpecial($xlPasteFormats)
$ws.Range("V235").Value = "Arduino, Tarjetas micro:bit."
$textFormatSource.Copy()
$ws.Range("W235").PasteSpecial($xlPasteFormats)
$ws.Range("W235").Value = 6.0
$textFormatSource.Copy()
$ws.Range("X235").PasteSpecial($xlPasteFormats)
$ws.Range("X235").Value = "No hay ningún elemento adicional a las micro:bits"
$textFormatSource.Copy()
$ws.Range("Z235").PasteSpecial($xlPasteFormats)
$ws.Range("Z235").Value = "No"
$textFormatSource.Copy()
$ws.Range("AA235").PasteSpecial($xlPasteFormats)
$ws.Range("AA235").Value = "No"
$textFormatSource.Copy()
$ws.Range("AB235").PasteSpecial($xlPasteFormats)
$ws.Range("AB235").Value = "No"

# ----- Row 236 -----
$dateFormatSource.Copy()
$ws.Range("A236").PasteSpecial($xlPasteFormats)
$ws.Range("A236").Value = 44698.302116990744
$textFormatSource.Copy()
$ws.Range("B236").PasteSpecial($xlPasteFormats)
$ws.Range("B236").Value = "a.nsf.fabian.rincon@cali.edu.co"
$textFormatSource.Copy()
$ws.Range("C236").PasteSpecial($xlPasteFormats)
$ws.Range("C236").Value = 17.0
$textFormatSource.Copy()
$ws.Range("D236").PasteSpecial($xlPasteFormats)
$ws.Range("D236").Value = "Una"
$textFormatSource.Copy()
$ws.Range("E236").PasteSpecial($xlPasteFormats)
$ws.Range("E236").Value = "Sí"
$textFormatSource.Copy()
$ws.Range("F236").PasteSpecial($xlPasteFormats)
$ws.Range("F236").Value = 16.0
$textFormatSource.Copy()
$ws.Range("G236").PasteSpecial($xlPasteFormats)
$ws.Range("G236").Value = "Video proyector"
$textFormatSource.Copy()
$ws.Range("H236").PasteSpecial($xlPasteFormats)
$ws.Range("H236").Value = "Puertos de cableado estructurado para conexión a la red"
$textFormatSource.Copy()
$ws.Range("J236").PasteSpecial($xlPasteFormats)
$ws.Range("J236").Value = 24.0
$textFormatSource.Copy()
$ws.Range("K236").PasteSpecial($xlPasteFormats)
$ws.Range("K236").Value = 24.0
$textFormatSource.Copy()
$ws.Range("L236").PasteSpecial($xlPasteFormats)
$ws.Range("L236").Value = 0.0
$textFormatSource.Copy()
$ws.Range("M236").PasteSpecial($xlPasteFormats)
$ws.Range("M236").Value = "Sí"
$textFormatSource.Copy()
$ws.Range("N236").PasteSpecial($xlPasteFormats)
$ws.Range("N236").Value = "Windows 7 o superior (32 bits)"
$textFormatSource.Copy()
$ws.Range("O236").PasteSpecial($xlPasteFormats)
$ws.Range("O236").Value = 0.0
$textFormatSource.Copy()
$ws.Range("P236").PasteSpecial($xlPasteFormats)
$ws.Range("P236").Value = "Paquete de oficina (Ej. Microsoft Office)"
$textFormatSource.Copy()
$ws.Range("Q236").PasteSpecial($xlPasteFormats)
$ws.Range("Q236").Value = "son portatiles mini "
$textFormatSource.Copy()
$ws.Range("R236").PasteSpecial($xlPasteFormats)
$ws.Range("R236").Value = "No"
$textFormatSource.Copy()
$ws.Range("V236").PasteSpecial($xlPasteFormats)
$ws.Range("V236").Value = "placa electronica y placa programable playboard"
$textFormatSource.Copy()
$ws.Range("W236").PasteSpecial($xlPasteFormats)
$ws.Range("W236").Value = 16.0
$textFormatSource.Copy()
$ws.Range("X236").PasteSpecial($xlPasteFormats)
$ws.Range("X236").Value = "Sensores, LEDs, Cables cocodrilo"
$textFormatSource.Copy()
$ws.Range("Y236").PasteSpecial($xlPasteFormats)
$ws.Range("Y236").Value = "sistema tecpro de norma que llego este año "
$textFormatSource.Copy()
$ws.Range("Z236").PasteSpecial($xlPasteFormats)
$ws.Range("Z236").Value = "Sí"
$textFormatSource.Copy()
$ws.Range("AA236").PasteSpecial($xlPasteFormats)
$ws.Range("AA236").Value = "Sí"
$textFormatSource.Copy()
$ws.Range("AB236").PasteSpecial($xlPasteFormats)
$ws.Range("AB236").Value = "No"
$textFormatSource.Copy()
$ws.Range("AC236").PasteSpecial($xlPasteFormats)
$ws.Range("AC236").Value = "se pueden reproducir hasta 30 ya que se cuenta con fotocopiadora"

# ----- Row 237 -----
$dateFormatSource.Copy()
$ws.Range("A237").PasteSpecial($xlPasteFormats)
$ws.Range("A237").Value = 44698.326642141205
$textFormatSource.Copy()
$ws.Range("B237").PasteSpecial($xlPasteFormats)
$ws.Range("B237").Value = "erika.miranda.cfk2022@gmail.com"
$textFormatSource.Copy()
$ws.Range("C237").PasteSpecial($xlPasteFormats)
$ws.Range("C237").Value = 85.0
$textFormatSource.Copy()
$ws.Range("D237").PasteSpecial($xlPasteFormats)
$ws.Range("D237").Value = "No tiene"
$textFormatSource.Copy()
$ws.Range("E237").PasteSpecial($xlPasteFormats)
$ws.Range("E237").Value = "No"
$textFormatSource.Copy()
$ws.Range("F237").PasteSpecial($xlPasteFormats)
$ws.Range("F237").Value = 0.0
$textFormatSource.Copy()
$ws.Range("G237").PasteSpecial($xlPasteFormats)
$ws.Range("G237").Value = "No hay más dispositivos en las aulas"
$textFormatSource.Copy()
$ws.Range("H237").PasteSpecial($xlPasteFormats)
$ws.Range("H237").Value = "Redes eléctricas, Puertos de cableado estructurado para conexión a la red"
$textFormatSource.Copy()
$ws.Range("I237").PasteSpecial($xlPasteFormats)
$ws.Range("I237").Value = "La IE tiene la infraestructura WIFI pero no hay servicio Internet."
$textFormatSource.Copy()
$ws.Range("J237").PasteSpecial($xlPasteFormats)
$ws.Range("J237").Value = 0.0
$textFormatSource.Copy()
$ws.Range("K237").PasteSpecial($xlPasteFormats)
$ws.Range("K237").Value = 0.0
$textFormatSource.Copy()
$ws.Range("L237").PasteSpecial($xlPasteFormats)
$ws.Range("L237").Value = 0.0
$textFormatSource.Copy()
$ws.Range("M237").PasteSpecial($xlPasteFormats)
$ws.Range("M237").Value = "No"
$textFormatSource.Copy()
$ws.Range("N237").PasteSpecial($xlPasteFormats)
$ws.Range("N237").Value = "No hay computadores"
$textFormatSource.Copy()
$ws.Range("O237").PasteSpecial($xlPasteFormats)
$ws.Range("O237").Value = 0.0
$textFormatSource.Copy()
$ws.Range("P237").PasteSpecial($xlPasteFormats)
$ws.Range("P237").Value = "Ninguno de los anteriores"
$textFormatSource.Copy()
$ws.Range("R237").PasteSpecial($xlPasteFormats)
$ws.Range("R237").Value = "No"
$textFormatSource.Copy()
$ws.Range("V237").PasteSpecial($xlPasteFormats)
$ws.Range("V237").Value = "Tarjetas micro:bit."
$textFormatSource.Copy()
$ws.Range("W237").PasteSpecial($xlPasteFormats)
$ws.Range("W237").Value = 12.0
$textFormatSource.Copy()
$ws.Range("X237").PasteSpecial($xlPasteFormats)
$ws.Range("X237").Value = "No hay ningún elemento adicional a las micro:bits"
$textFormatSource.Copy()
$ws.Range("Z237").PasteSpecial($xlPasteFormats)
$ws.Range("Z237").Value = "No"
$textFormatSource.Copy()
$ws.Range("AA237").PasteSpecial($xlPasteFormats)
$ws.Range("AA237").Value = "No"
$textFormatSource.Copy()
$ws.Range("AB237").PasteSpecial($xlPasteFormats)
$ws.Range("AB237").Value = "No"

# ----- Row 238 -----
$dateFormatSource.Copy()
$ws.Range("A238").PasteSpecial($xlPasteFormats)
$ws.Range("A238").Value = 44698.464638958336
$textFormatSource.Copy()
$ws.Range("B238").PasteSpecial($xlPasteFormats)
$ws.Range("B238").Value = "carolina.timana.cfk2022@gmail.com"
$textFormatSource.Copy()
$ws.Range("C238").PasteSpecial($xlPasteFormats)
$ws.Range("C238").Value = 135.0
$textFormatSource.Copy()
$ws.Range("D238").PasteSpecial($xlPasteFormats)
$ws.Range("D238").Value = "Dos"
$textFormatSource.Copy()
$ws.Range("E238").PasteSpecial($xlPasteFormats)
$ws.Range("E238").Value = "No"
$textFormatSource.Copy()
$ws.Range("F238").PasteSpecial($xlPasteFormats)
$ws.Range("F238").Value = 0.0
$textFormatSource.Copy()
$ws.Range("G238").PasteSpecial($xlPasteFormats)
$ws.Range("G238").Value = "No hay más dispositivos en las aulas"
$textFormatSource.Copy()
$ws.Range("H238").PasteSpecial($xlPasteFormats)
$ws.Range("H238").Value = "Redes eléctricas"
$textFormatSource.Copy()
$ws.Range("J238").PasteSpecial($xlPasteFormats)
$ws.Range("J238").Value = 37.0
$textFormatSource.Copy()
$ws.Range("K238").PasteSpecial($xlPasteFormats)
$ws.Range("K238").Value = 25.0
$textFormatSource.Copy()
$ws.Range("L238").PasteSpecial($xlPasteFormats)
$ws.Range("L238").Value = 0.0
$textFormatSource.Copy()
$ws.Range("M238").PasteSpecial($xlPasteFormats)
$ws.Range("M238").Value = "No"
$textFormatSource.Copy()
$ws.Range("N238").PasteSpecial($xlPasteFormats)
$ws.Range("N238").Value = "Windows 7 o superior (64 bits)"
$textFormatSource.Copy()
$ws.Range("O238").PasteSpecial($xlPasteFormats)
$ws.Range("O238").Value = 0.0
$textFormatSource.Copy()
$ws.Range("P238").PasteSpecial($xlPasteFormats)
$ws.Range("P238").Value = "Paquete de oficina (Ej. Microsoft Office), Paquete de software educativo (Ej. cursos de inglés)"
$textFormatSource.Copy()
$ws.Range("Q238").PasteSpecial($xlPasteFormats)
$ws.Range("Q238").Value = "Se requiere actualizar los equipos e instalar nuevas aulas de sistemas"
$textFormatSource.Copy()
$ws.Range("R238").PasteSpecial($xlPasteFormats)
$ws.Range("R238").Value = "No"
$textFormatSource.Copy()
$ws.Range("V238").PasteSpecial($xlPasteFormats)
$ws.Range("V238").Value = "Tarjetas micro:bit."
$textFormatSource.Copy()
$ws.Range("W238").PasteSpecial($xlPasteFormats)
$ws.Range("W238").Value = 6.0
$textFormatSource.Copy()
$ws.Range("X238").PasteSpecial($xlPasteFormats)
$ws.Range("X238").Value = "No hay ningún elemento adicional a las micro:bits"
$textFormatSource.Copy()
$ws.Range("Y238").PasteSpecial($xlPasteFormats)
$ws.Range("Y238").Value = "Las tabletas y los equipos portátiles de computadores para educar fueron robados en el tiempo de la pandemia"
$textFormatSource.Copy()
$ws.Range("Z238").PasteSpecial($xlPasteFormats)
$ws.Range("Z238").Value = "No"
$textFormatSource.Copy()
$ws.Range("AA238").PasteSpecial($xlPasteFormats)
$ws.Range("AA238").Value = "No"
$textFormatSource.Copy()
$ws.Range("AB238").PasteSpecial($xlPasteFormats)
$ws.Range("AB238").Value = "No"
$textFormatSource.Copy()
$ws.Range("AC238").PasteSpecial($xlPasteFormats)
$ws.Range("AC238").Value = "Existen 6 maletas Handy cricket pero solo hay las fichas lego ya no existen los componentes electrónicos"

# ----- Row 239 -----
$dateFormatSource.Copy()
$ws.Range("A239").PasteSpecial($xlPasteFormats)
$ws.Range("A239").Value = 44698.627339513885
$textFormatSource.Copy()
$ws.Range("B239").PasteSpecial($xlPasteFormats)
$ws.Range("B239").Value = "monica.giraldo.cfk2022@gmail.com"
$textFormatSource.Copy()
$ws.Range("C239").PasteSpecial($xlPasteFormats)
$ws.Range("C239").Value = 233.0
$textFormatSource.Copy()
$ws.Range("D239").PasteSpecial($xlPasteFormats)
$ws.Range("D239").Value = "Tres"
$textFormatSource.Copy()
$ws.Range("E239").PasteSpecial($xlPasteFormats)
$ws.Range("E239").Value = "Algunas sí, pero no todas"
$textFormatSource.Copy()
$ws.Range("F239").PasteSpecial($xlPasteFormats)
$ws.Range("F239").Value = 3.0
$textFormatSource.Copy()
$ws.Range("G239").PasteSpecial($xlPasteFormats)
$ws.Range("G239").Value = "Televisor"
$textFormatSource.Copy()
$ws.Range("H239").PasteSpecial($xlPasteFormats)
$ws.Range("H239").Value = "Redes eléctricas, Puertos de cableado estructurado para conexión a la red"
$textFormatSource.Copy()
$ws.Range("I239").PasteSpecial($xlPasteFormats)
$ws.Range("I239").Value = "En pocas aulas hay televisor"
$textFormatSource.Copy()
$ws.Range("J239").PasteSpecial($xlPasteFormats)
$ws.Range("J239").Value = 71.0
$textFormatSource.Copy()
$ws.Range("K239").PasteSpecial($xlPasteFormats)
$ws.Range("K239").Value = 31.0
$textFormatSource.Copy()
$ws.Range("L239").PasteSpecial($xlPasteFormats)
$ws.Range("L239").Value = 0.0
$textFormatSource.Copy()
$ws.Range("M239").PasteSpecial($xlPasteFormats)
$ws.Range("M239").Value = "Algunos de ellos"
$textFormatSource.Copy()
$ws.Range("N239").PasteSpecial($xlPasteFormats)
$ws.Range("N239").Value = "Windows 7 o superior (32 bits)"
$textFormatSource.Copy()
$ws.Range("O239").PasteSpecial($xlPasteFormats)
$ws.Range("O239").Value = 0.0
$textFormatSource.Copy()
$ws.Range("P239").PasteSpecial($xlPasteFormats)
$ws.Range("P239").Value = "Paquete de oficina (Ej. Microsoft Office)"
$textFormatSource.Copy()
$ws.Range("Q239").PasteSpecial($xlPasteFormats)
$ws.Range("Q239").Value = "El paquete de oficina se encuentra instalado sólo en algunos computadores. Los que se conectan en la nube no permiten ningún tipo de descarga y limita accesos "
$textFormatSource.Copy()
$ws.Range("R239").PasteSpecial($xlPasteFormats)
$ws.Range("R239").Value = "No"
$textFormatSource.Copy()
$ws.Range("V239").PasteSpecial($xlPasteFormats)
$ws.Range("V239").Value = "Tarjetas micro:bit."
$textFormatSource.Copy()
$ws.Range("W239").PasteSpecial($xlPasteFormats)
$ws.Range("W239").Value = 28.0
$textFormatSource.Copy()
$ws.Range("X239").PasteSpecial($xlPasteFormats)
$ws.Range("X239").Value = "No hay ningún elemento adicional a las micro:bits"
$textFormatSource.Copy()
$ws.Range("Y239").PasteSpecial($xlPasteFormats)
$ws.Range("Y239").Value = "La IE cuenta con unos kit de robótica, sin embargo no se logró acceso a ellos para determinar que contienen y su estado. "
$textFormatSource.Copy()
$ws.Range("Z239").PasteSpecial($xlPasteFormats)
$ws.Range("Z239").Value = "No"
$textFormatSource.Copy()
$ws.Range("AA239").PasteSpecial($xlPasteFormats)
$ws.Range("AA239").Value = "No"
$textFormatSource.Copy()
$ws.Range("AB239").PasteSpecial($xlPasteFormats)
$ws.Range("AB239").Value = "No"

# ----- Row 240 -----
$dateFormatSource.Copy()
$ws.Range("A240").PasteSpecial($xlPasteFormats)
$ws.Range("A240").Value = 44699.415555451385
$textFormatSource.Copy()
$ws.Range("B240").PasteSpecial($xlPasteFormats)
$ws.Range("B240").Value = "camilo.villota.cfk2022@gmail.com"
$textFormatSource.Copy()
$ws.Range("C240").PasteSpecial($xlPasteFormats)
$ws.Range("C240").Value = 38.0
$textFormatSource.Copy()
$ws.Range("D240").PasteSpecial($xlPasteFormats)
$ws.Range("D240").Value = "Dos"
$textFormatSource.Copy()
$ws.Range("E240").PasteSpecial($xlPasteFormats)
$ws.Range("E240").Value = "No"
$textFormatSource.Copy()
$ws.Range("F240").PasteSpecial($xlPasteFormats)
$ws.Range("F240").Value = 0.0
$textFormatSource.Copy()
$ws.Range("G240").PasteSpecial($xlPasteFormats)
$ws.Range("G240").Value = "Televisor, Video proyector, Pero no están en la sala, se tienen pero se usan cuando los profesores, 11 televisores, 5 videobeam."
$textFormatSource.Copy()
$ws.Range("H240").PasteSpecial($xlPasteFormats)
$ws.Range("H240").Value = "Redes eléctricas, Puertos de cableado estructurado para conexión a la red, Convenio alcaldía-ETB para el tema de conectividad pero no funciona."
$textFormatSource.Copy()
$ws.Range("J240").PasteSpecial($xlPasteFormats)
$ws.Range("J240").Value = 220.0
$textFormatSource.Copy()
$ws.Range("K240").PasteSpecial($xlPasteFormats)
$ws.Range("K240").Value = 184.0
$textFormatSource.Copy()
$ws.Range("L240").PasteSpecial($xlPasteFormats)
$ws.Range("L240").Value = 0.0
$textFormatSource.Copy()
$ws.Range("M240").PasteSpecial($xlPasteFormats)
$ws.Range("M240").Value = "Algunos de ellos"
$textFormatSource.Copy()
$ws.Range("N240").PasteSpecial($xlPasteFormats)
$ws.Range("N240").Value = "Windows 7 o superior (64 bits)"
$textFormatSource.Copy()
$ws.Range("O240").PasteSpecial($xlPasteFormats)
$ws.Range("O240").Value = 0.0
$textFormatSource.Copy()
$ws.Range("P240").PasteSpecial($xlPasteFormats)
$ws.Range("P240").Value = "Paquete de oficina (Ej. Microsoft Office)"
$textFormatSource.Copy()
$ws.Range("R240").PasteSpecial($xlPasteFormats)
$ws.Range("R240").Value = "Sí"
$textFormatSource.Copy()
$ws.Range("S240").PasteSpecial($xlPasteFormats)
$ws.Range("S240").Value = 220.0
$textFormatSource.Copy()
$ws.Range("T240").PasteSpecial($xlPasteFormats)
$ws.Range("T240").Value = "Sí"
$textFormatSource.Copy()
$ws.Range("U240").PasteSpecial($xlPasteFormats)
$ws.Range("U240").Value = 0.0
$textFormatSource.Copy()
$ws.Range("V240").PasteSpecial($xlPasteFormats)
$ws.Range("V240").Value = "Afirman no tener las micro:bit"
$textFormatSource.Copy()
$ws.Range("W240").PasteSpecial($xlPasteFormats)
$ws.Range("W240").Value = 0.0
$textFormatSource.Copy()
$ws.Range("X240").PasteSpecial($xlPasteFormats)
$ws.Range("X240").Value = "No hay ningún elemento adicional a las micro:bits"
$textFormatSource.Copy()
$ws.Range("Y240").PasteSpecial($xlPasteFormats)
$ws.Range("Y240").Value = "Al parecer la institución no ha recibido las micro:bit que se entragaron a docentes porque no están en inventario."
$textFormatSource.Copy()
$ws.Range("Z240").PasteSpecial($xlPasteFormats)
$ws.Range("Z240").Value = "No"
$textFormatSource.Copy()
$ws.Range("AA240").PasteSpecial($xlPasteFormats)
$ws.Range("AA240").Value = "No"
$textFormatSource.Copy()
$ws.Range("AB240").PasteSpecial($xlPasteFormats)
$ws.Range("AB240").Value = "No"

# ----- Row 241 -----
$dateFormatSource.Copy()
$ws.Range("A241").PasteSpecial($xlPasteFormats)
$ws.Range("A241").Value = 44699.435526967594
$textFormatSource.Copy()
$ws.Range("B241").PasteSpecial($xlPasteFormats)
$ws.Range("B241").Value = "jhon.balcarcel.cfk2022@gmail.com"
$textFormatSource.Copy()
$ws.Range("C241").PasteSpecial($xlPasteFormats)
$ws.Range("C241").Value = 146.0
$textFormatSource.Copy()
$ws.Range("D241").PasteSpecial($xlPasteFormats)
$ws.Range("D241").Value = "Dos"
$textFormatSource.Copy()
$ws.Range("E241").PasteSpecial($xlPasteFormats)
$ws.Range("E241").Value = "No"
$textFormatSource.Copy()
$ws.Range("F241").PasteSpecial($xlPasteFormats)
$ws.Range("F241").Value = 0.0
$textFormatSource.Copy()
$ws.Range("G241").PasteSpecial($xlPasteFormats)
$ws.Range("G241").Value = "Televisor"
$textFormatSource.Copy()
$ws.Range("H241").PasteSpecial($xlPasteFormats)
$ws.Range("H241").Value = "Redes eléctricas, Puertos de cableado estructurado para conexión a la red, Módem"
$textFormatSource.Copy()
$ws.Range("I241").PasteSpecial($xlPasteFormats)
$ws.Range("I241").Value = "Un aula de sistemas esta en proceso de adecuación y hay otra funcional"
$textFormatSource.Copy()
$ws.Range("J241").PasteSpecial($xlPasteFormats)
$ws.Range("J241").Value = 20.0
$textFormatSource.Copy()
$ws.Range("K241").PasteSpecial($xlPasteFormats)
$ws.Range("K241").Value = 20.0
$textFormatSource.Copy()
$ws.Range("L241").PasteSpecial($xlPasteFormats)
$ws.Range("L241").Value = 0.0
$textFormatSource.Copy()
$ws.Range("M241").PasteSpecial($xlPasteFormats)
$ws.Range("M241").Value = "Algunos de ellos"
$textFormatSource.Copy()
$ws.Range("N241").PasteSpecial($xlPasteFormats)
$ws.Range("N241").Value = "Windows 7 o superior (32 bits)"
$textFormatSource.Copy()
$ws.Range("O241").PasteSpecial($xlPasteFormats)
$ws.Range("O241").Value = 0.0
$textFormatSource.Copy()
$ws.Range("P241").PasteSpecial($xlPasteFormats)
$ws.Range("P241").Value = "Paquete de oficina (Ej. Microsoft Office)"
$textFormatSource.Copy()
$ws.Range("R241").PasteSpecial($xlPasteFormats)
$ws.Range("R241").Value = "No"
$textFormatSource.Copy()
$ws.Range("V241").PasteSpecial($xlPasteFormats)
$ws.Range("V241").Value = "Tarjetas micro:bit."
$textFormatSource.Copy()
$ws.Range("W241").PasteSpecial($xlPasteFormats)
$ws.Range("W241").Value = 12.0
$textFormatSource.Copy()
$ws.Range("X241").PasteSpecial($xlPasteFormats)
$ws.Range("X241").Value = "No hay ningún elemento adicional a las micro:bits"
$textFormatSource.Copy()
$ws.Range("Y241").PasteSpecial($xlPasteFormats)
$ws.Range("Y241").Value = "Sòlo estan limitados a los equipos de computo y las microbit"
$textFormatSource.Copy()
$ws.Range("Z241").PasteSpecial($xlPasteFormats)
$ws.Range("Z241").Value = "No"
$textFormatSource.Copy()
$ws.Range("AA241").PasteSpecial($xlPasteFormats)
$ws.Range("AA241").Value = "No"
$textFormatSource.Copy()
$ws.Range("AB241").PasteSpecial($xlPasteFormats)
$ws.Range("AB241").Value = "No"

# ----- Row 242 -----
$dateFormatSource.Copy()
$ws.Range("A242").PasteSpecial($xlPasteFormats)
$ws.Range("A242").Value = 44699.45864196759
$textFormatSource.Copy()
$ws.Range("B242").PasteSpecial($xlPasteFormats)
$ws.Range("B242").Value = "yessicaeverts@gmail.com"
$textFormatSource.Copy()
$ws.Range("C242").PasteSpecial($xlPasteFormats)
$ws.Range("C242").Value = 194.0
$textFormatSource.Copy()
$ws.Range("D242").PasteSpecial($xlPasteFormats)
$ws.Range("D242").Value = "Una"
$textFormatSource.Copy()
$ws.Range("E242").PasteSpecial($xlPasteFormats)
$ws.Range("E242").Value = "No"
$textFormatSource.Copy()
$ws.Range("F242").PasteSpecial($xlPasteFormats)
$ws.Range("F242").Value = 0.0
$textFormatSource.Copy()
$ws.Range("G242").PasteSpecial($xlPasteFormats)
$ws.Range("G242").Value = "Televisor, Video proyector"
$textFormatSource.Copy()
$ws.Range("H242").PasteSpecial($xlPasteFormats)
$ws.Range("H242").Value = "No hay ninguna infraestructura de telecomunicaciones"
$textFormatSource.Copy()
$ws.Range("I242").PasteSpecial($xlPasteFormats)
$ws.Range("I242").Value = "Se acabó el contrato de conexión a internet en el colegio y no se ha renovado. "
$textFormatSource.Copy()
$ws.Range("J242").PasteSpecial($xlPasteFormats)
$ws.Range("J242").Value = 41.0
$textFormatSource.Copy()
$ws.Range("K242").PasteSpecial($xlPasteFormats)
$ws.Range("K242").Value = 14.0
$textFormatSource.Copy()
$ws.Range("L242").PasteSpecial($xlPasteFormats)
$ws.Range("L242").Value = 0.0
$textFormatSource.Copy()
$ws.Range("M242").PasteSpecial($xlPasteFormats)
$ws.Range("M242").Value = "No"
$textFormatSource.Copy()
$ws.Range("N242").PasteSpecial($xlPasteFormats)
$ws.Range("N242").Value = "Windows 7 o superior (64 bits)"
$textFormatSource.Copy()
$ws.Range("O242").PasteSpecial($xlPasteFormats)
$ws.Range("O242").Value = 0.0
$textFormatSource.Copy()
$ws.Range("P242").PasteSpecial($xlPasteFormats)
$ws.Range("P242").Value = "Paquete de oficina (Ej. Microsoft Office)"
$textFormatSource.Copy()
$ws.Range("R242").PasteSpecial($xlPasteFormats)
$ws.Range("R242").Value = "Sí"
$textFormatSource.Copy()
$ws.Range("S242").PasteSpecial($xlPasteFormats)
$ws.Range("S242").Value = 50.0
$textFormatSource.Copy()
$ws.Range("T242").PasteSpecial($xlPasteFormats)
$ws.Range("T242").Value = "Sí"
$textFormatSource.Copy()
$ws.Range("U242").PasteSpecial($xlPasteFormats)
$ws.Range("U242").Value = 0.0
$textFormatSource.Copy()
$ws.Range("V242").PasteSpecial($xlPasteFormats)
$ws.Range("V242").Value = "Tarjetas micro:bit."
$textFormatSource.Copy()
$ws.Range("W242").PasteSpecial($xlPasteFormats)
$ws.Range("W242").Value = 1.0
$textFormatSource.Copy()
$ws.Range("X242").PasteSpecial($xlPasteFormats)
$ws.Range("X242").Value = "No hay ningún elemento adicional a las micro:bits"
$textFormatSource.Copy()
$ws.Range("Z242").PasteSpecial($xlPasteFormats)
$ws.Range("Z242").Value = "No"
$textFormatSource.Copy()
$ws.Range("AA242").PasteSpecial($xlPasteFormats)
$ws.Range("AA242").Value = "No"
$textFormatSource.Copy()
$ws.Range("AB242").PasteSpecial($xlPasteFormats)
$ws.Range("AB242").Value = "No"

# ----- Row 243 -----
$dateFormatSource.Copy()
$ws.Range("A243").PasteSpecial($xlPasteFormats)
$ws.Range("A243").Value = 44699.99359429398
$textFormatSource.Copy()
$ws.Range("B243").PasteSpecial($xlPasteFormats)
$ws.Range("B243").Value = "sami123rojare@yahoo.es"
$textFormatSource.Copy()
$ws.Range("C243").PasteSpecial($xlPasteFormats)
$ws.Range("C243").Value = 207.0
$textFormatSource.Copy()
$ws.Range("D243").PasteSpecial($xlPasteFormats)
$ws.Range("D243").Value = "Una"
$textFormatSource.Copy()
$ws.Range("E243").PasteSpecial($xlPasteFormats)
$ws.Range("E243").Value = "No"
$textFormatSource.Copy()
$ws.Range("F243").PasteSpecial($xlPasteFormats)
$ws.Range("F243").Value = 0.0
$textFormatSource.Copy()
$ws.Range("G243").PasteSpecial($xlPasteFormats)
$ws.Range("G243").Value = "Televisor, Video proyector"
$textFormatSource.Copy()
$ws.Range("H243").PasteSpecial($xlPasteFormats)
$ws.Range("H243").Value = "Redes eléctricas, Módem"
$textFormatSource.Copy()
$ws.Range("I243").PasteSpecial($xlPasteFormats)
$ws.Range("I243").Value = "En la IE, se presenta mala señal de celular y la conectividad es muy deficiente para garantizar conexion a los equipos de computo de la sala."
$textFormatSource.Copy()
$ws.Range("J243").PasteSpecial($xlPasteFormats)
$ws.Range("J243").Value = 40.0
$textFormatSource.Copy()
$ws.Range("K243").PasteSpecial($xlPasteFormats)
$ws.Range("K243").Value = 30.0
$textFormatSource.Copy()
$ws.Range("L243").PasteSpecial($xlPasteFormats)
$ws.Range("L243").Value = 0.0
$textFormatSource.Copy()
$ws.Range("M243").PasteSpecial($xlPasteFormats)
$ws.Range("M243").Value = "Algunos de ellos"
$textFormatSource.Copy()
$ws.Range("N243").PasteSpecial($xlPasteFormats)
$ws.Range("N243").Value = "Windows 7 o superior (64 bits)"
$textFormatSource.Copy()
$ws.Range("O243").PasteSpecial($xlPasteFormats)
$ws.Range("O243").Value = 0.0
$textFormatSource.Copy()
$ws.Range("P243").PasteSpecial($xlPasteFormats)
$ws.Range("P243").Value = "Paquete de diseño gráfico (Ej. Programas de edición de video, edición de mapas de bits, etc), Paquete de software educativo (Ej. cursos de inglés), no se tiene office instalado en los equipos"
$textFormatSource.Copy()
$ws.Range("Q243").PasteSpecial($xlPasteFormats)
$ws.Range("Q243").Value = "actualmente el docente de tecnologia, esta investigando sobre que programas utilizar en el aula de clase."
$textFormatSource.Copy()
$ws.Range("R243").PasteSpecial($xlPasteFormats)
$ws.Range("R243").Value = "No"
$textFormatSource.Copy()
$ws.Range("V243").PasteSpecial($xlPasteFormats)
$ws.Range("V243").Value = "Tarjetas micro:bit., solo 3 microbits"
$textFormatSource.Copy()
$ws.Range("W243").PasteSpecial($xlPasteFormats)
$ws.Range("W243").Value = 3.0
$textFormatSource.Copy()
$ws.Range("X243").PasteSpecial($xlPasteFormats)
$ws.Range("X243").Value = "No hay ningún elemento adicional a las micro:bits"
$textFormatSource.Copy()
$ws.Range("Y243").PasteSpecial($xlPasteFormats)
$ws.Range("Y243").Value = "se cuenta con una buena sala de tecnología, pero carece de buena conectividad, en la ie recién instalaron un centro poblado de conectividad pero nunca funciona y no suple la capacidad del numero de estudiantes."
$textFormatSource.Copy()
$ws.Range("Z243").PasteSpecial($xlPasteFormats)
$ws.Range("Z243").Value = "No"
$textFormatSource.Copy()
$ws.Range("AA243").PasteSpecial($xlPasteFormats)
$ws.Range("AA243").Value = "No"
$textFormatSource.Copy()
$ws.Range("AB243").PasteSpecial($xlPasteFormats)
$ws.Range("AB243").Value = "Sí"
$textFormatSource.Copy()
$ws.Range("AC243").PasteSpecial($xlPasteFormats)
$ws.Range("AC243").Value = "solo se encuentran las cartillas que la docente formada ha realizado"

# ----- Row 244 -----
$dateFormatSource.Copy()
$ws.Range("A244").PasteSpecial($xlPasteFormats)
$ws.Range("A244").Value = 44700.37422266204
$textFormatSource.Copy()
$ws.Range("B244").PasteSpecial($xlPasteFormats)
$ws.Range("B244").Value = "carolina.timana.cfk2022@gmail.com"
$textFormatSource.Copy()
$ws.Range("C244").PasteSpecial($xlPasteFormats)
$ws.Range("C244").Value = 137.0
$textFormatSource.Copy()
$ws.Range("D244").PasteSpecial($xlPasteFormats)
$ws.Range("D244").Value = "Dos"
$textFormatSource.Copy()
$ws.Range("E244").PasteSpecial($xlPasteFormats)
$ws.Range("E244").Value = "No"
$textFormatSource.Copy()
$ws.Range("F244").PasteSpecial($xlPasteFormats)
$ws.Range("F244").Value = 0.0
$textFormatSource.Copy()
$ws.Range("G244").PasteSpecial($xlPasteFormats)
$ws.Range("G244").Value = "No hay más dispositivos en las aulas"
$textFormatSource.Copy()
$ws.Range("H244").PasteSpecial($xlPasteFormats)
$ws.Range("H244").Value = "Redes eléctricas, Puertos de cableado estructurado para conexión a la red, Módem"
$textFormatSource.Copy()
$ws.Range("J244").PasteSpecial($xlPasteFormats)
$ws.Range("J244").Value = 90.0
$textFormatSource.Copy()
$ws.Range("K244").PasteSpecial($xlPasteFormats)
$ws.Range("K244").Value = 88.0
$textFormatSource.Copy()
$ws.Range("L244").PasteSpecial($xlPasteFormats)
$ws.Range("L244").Value = 0.0
$textFormatSource.Copy()
$ws.Range("M244").PasteSpecial($xlPasteFormats)
$ws.Range("M244").Value = "Algunos de ellos"
$textFormatSource.Copy()
$ws.Range("N244").PasteSpecial($xlPasteFormats)
$ws.Range("N244").Value = "Windows 7 o superior (64 bits)"
$textFormatSource.Copy()
$ws.Range("O244").PasteSpecial($xlPasteFormats)
$ws.Range("O244").Value = 0.0
$textFormatSource.Copy()
$ws.Range("P244").PasteSpecial($xlPasteFormats)
$ws.Range("P244").Value = "Paquete de oficina (Ej. Microsoft Office)"
$textFormatSource.Copy()
$ws.Range("R244").PasteSpecial($xlPasteFormats)
$ws.Range("R244").Value = "Sí"
$textFormatSource.Copy()
$ws.Range("S244").PasteSpecial($xlPasteFormats)
$ws.Range("S244").Value = 24.0
$textFormatSource.Copy()
$ws.Range("T244").PasteSpecial($xlPasteFormats)
$ws.Range("T244").Value = "No"
$textFormatSource.Copy()
$ws.Range("U244").PasteSpecial($xlPasteFormats)
$ws.Range("U244").Value = 0.0
$textFormatSource.Copy()
$ws.Range("V244").PasteSpecial($xlPasteFormats)
$ws.Range("V244").Value = "Arduino, Tarjetas micro:bit."
$textFormatSource.Copy()
$ws.Range("W244").PasteSpecial($xlPasteFormats)
$ws.Range("W244").Value = 9.0
$textFormatSource.Copy()
$ws.Range("X244").PasteSpecial($xlPasteFormats)
$ws.Range("X244").Value = "LEDs, Servomotores, Cables cocodrilo"
$textFormatSource.Copy()
$ws.Range("Z244").PasteSpecial($xlPasteFormats)
$ws.Range("Z244").Value = "No"
$textFormatSource.Copy()
$ws.Range("AA244").PasteSpecial($xlPasteFormats)
$ws.Range("AA244").Value = "Sí"
$textFormatSource.Copy()
$ws.Range("AB244").PasteSpecial($xlPasteFormats)
$ws.Range("AB244").Value = "No"
$textFormatSource.Copy()
$ws.Range("AC244").PasteSpecial($xlPasteFormats)
$ws.Range("AC244").Value = "Existe en la Institución un (1) kit de ruta STEM"
